$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 496-503 (values change in place, columns A,B,C,E-K,R constant) ---
$ws.Range("D496").Value2 = 44595
$ws.Range("L496").Value2 = "Especial"
$ws.Range("M496").Value2 = 162
$ws.Range("N496").Value2 = 17000
$ws.Range("O496").Value2 = 17000
$ws.Range("P496").Value2 = 17000
$ws.Range("Q496").Value2 = "`$/caja 10 unidades"
$ws.Range("S496").Value2 = 1700
$ws.Range("T496").Value2 = 10

$ws.Range("D497").Value2 = 44595
$ws.Range("L497").Value2 = "Primera"
$ws.Range("M497").Value2 = 162
$ws.Range("N497").Value2 = 17000
$ws.Range("O497").Value2 = 17000
$ws.Range("P497").Value2 = 17000
$ws.Range("Q497").Value2 = "`$/caja 12 unidades"
$ws.Range("S497").Value2 = 1417
$ws.Range("T497").Value2 = 12

$ws.Range("D498").Value2 = 44595
$ws.Range("L498").Value2 = "Segunda"
$ws.Range("M498").Value2 = 108
$ws.Range("N498").Value2 = 17000
$ws.Range("O498").Value2 = 17000
$ws.Range("P498").Value2 = 17000
$ws.Range("Q498").Value2 = "`$/caja 14 unidades"
$ws.Range("S498").Value2 = 1214
$ws.Range("T498").Value2 = 14

$ws.Range("D499").Value2 = 44335
$ws.Range("L499").Value2 = "Primera"
$ws.Range("M499").Value2 = 162
$ws.Range("N499").Value2 = 17000
$ws.Range("O499").Value2 = 17000
$ws.Range("P499").Value2 = 17000
$ws.Range("Q499").Value2 = "`$/caja 12 unidades"
$ws.Range("S499").Value2 = 1417
$ws.Range("T499").Value2 = 12

$ws.Range("D500").Value2 = 44552
$ws.Range("L500").Value2 = "Especial"
$ws.Range("M500").Value2 = 108
$ws.Range("N500").Value2 = 18000
$ws.Range("O500").Value2 = 18000
$ws.Range("P500").Value2 = 18000
$ws.Range("Q500").Value2 = "`$/caja 10 unidades"
$ws.Range("S500").Value2 = 1800
$ws.Range("T500").Value2 = 10

$ws.Range("D501").Value2 = 44552
$ws.Range("L501").Value2 = "Primera"
$ws.Range("M501").Value2 = 162
$ws.Range("N501").Value2 = 18000
$ws.Range("O501").Value2 = 18000
$ws.Range("P501").Value2 = 18000
$ws.Range("Q501").Value2 = "`$/caja 12 unidades"
$ws.Range("S501").Value2 = 1500
$ws.Range("T501").Value2 = 12

$ws.Range("D502").Value2 = 44552
$ws.Range("L502").Value2 = "Segunda"
$ws.Range("M502").Value2 = 108
$ws.Range("N502").Value2 = 18000
$ws.Range("O502").Value2 = 18000
$ws.Range("P502").Value2 = 18000
$ws.Range("Q502").Value2 = "`$/caja 14 unidades"
$ws.Range("S502").Value2 = 1286
$ws.Range("T502").Value2 = 14

$ws.Range("D503").Value2 = 44544
$ws.Range("L503").Value2 = "Especial"
$ws.Range("M503").Value2 = 108
$ws.Range("N503").Value2 = 18000
$ws.Range("O503").Value2 = 18000
$ws.Range("P503").Value2 = 18000
$ws.Range("Q503").Value2 = "`$/caja 10 unidades"
$ws.Range("S503").Value2 = 1800
$ws.Range("T503").Value2 = 10

# --- Row 504 becomes new content (previously row 504 was "Segunda/22000", now replaced by shifted "Primera/18000" data) ---
$ws.Range("D504").Value2 = 44544
$ws.Range("L504").Value2 = "Primera"
$ws.Range("M504").Value2 = 162
$ws.Range("N504").Value2 = 18000
$ws.Range("O504").Value2 = 18000
$ws.Range("P504").Value2 = 18000
$ws.Range("Q504").Value2 = "`$/caja 12 unidades"
$ws.Range("S504").Value2 = 1500
$ws.Range("T504").Value2 = 12

# --- New rows 505-507: constant columns copied, then variable columns set ---
$ws.Range("A505").Value2 = 3
$ws.Range("B505").Value2 = "Femacal de La Calera"
$ws.Range("C505").Value2 = "Coquimbo"
$ws.Range("D505").Value2 = 44544
$ws.Range("D505").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E505").Value2 = 5
$ws.Range("F505").Value2 = "Fruta"
$ws.Range("G505").Value2 = 100108
$ws.Range("H505").Value2 = "Tropicales y subtropicales"
$ws.Range("I505").Value2 = 100108005
$ws.Range("J505").Value2 = "Piña"
$ws.Range("K505").Value2 = "Caramelo"
$ws.Range("L505").Value2 = "Segunda"
$ws.Range("M505").Value2 = 108
$ws.Range("N505").Value2 = 18000
$ws.Range("O505").Value2 = 18000
$ws.Range("P505").Value2 = 18000
$ws.Range("Q505").Value2 = "`$/caja 14 unidades"
$ws.Range("R505").Value2 = "Ecuador"
$ws.Range("S505").Value2 = 1286
$ws.Range("T505").Value2 = 14

$ws.Range("A506").Value2 = 3
$ws.Range("B506").Value2 = "Femacal de La Calera"
$ws.Range("C506").Value2 = "Coquimbo"
$ws.Range("D506").Value2 = 44160
$ws.Range("D506").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E506").Value2 = 5
$ws.Range("F506").Value2 = "Fruta"
$ws.Range("G506").Value2 = 100108
$ws.Range("H506").Value2 = "Tropicales y subtropicales"
$ws.Range("I506").Value2 = 100108005
$ws.Range("J506").Value2 = "Piña"
$ws.Range("K506").Value2 = "Caramelo"
$ws.Range("L506").Value2 = "Primera"
$ws.Range("M506").Value2 = 162
$ws.Range("N506").Value2 = 22000
$ws.Range("O506").Value2 = 22000
$ws.Range("P506").Value2 = 22000
$ws.Range("Q506").Value2 = "`$/caja 12 unidades"
$ws.Range("R506").Value2 = "Ecuador"
$ws.Range("S506").Value2 = 1833
$ws.Range("T506").Value2 = 12

$ws.Range("A507").Value2 = 3
$ws.Range("B507").Value2 = "Femacal de La Calera"
$ws.Range("C507").Value2 = "Coquimbo"
$ws.Range("D507").Value2 = 44160
$ws.Range("D507").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E507").Value2 = 5
$ws.Range("F507").Value2 = "Fruta"
$ws.Range("G507").Value2 = 100108
$ws.Range("H507").Value2 = "Tropicales y subtropicales"
$ws.Range("I507").Value2 = 100108005
$ws.Range("J507").Value2 = "Piña"
$ws.Range("K507").Value2 = "Caramelo"
$ws.Range("L507").Value2 = "Segunda"
$ws.Range("M507").Value2 = 108
$ws.Range("N507").Value2 = 22000
$ws.Range("O507").Value2 = 22000
$ws.Range("P507").Value2 = 22000
$ws.Range("Q507").Value2 = "`$/caja 14 unidades"
$ws.Range("R507").Value2 = "Ecuador"
$ws.Range("S507").Value2 = 1571
$ws.Range("T507").Value2 = 14
